$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "add"
$ws.Range("C2").Value = "new"
$ws.Range("D2").Value = "line"

$ws.Range("G7").Select()
